# Generate Report for Handback
# Update timestamp values that reflect a refreshed handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for fe721b5d-0fa7-434f-8789-e65e1e9fafb1.md
$wsOverview.Range("G4").Value = "2016-08-21 00:52:15"

# zh-cn sheet: Correspond Handoff / Handback DateTime for fe721b5d-0fa7-434f-8789-e65e1e9fafb1 row
$wsZhCn.Range("H4").Value = "2016-08-21 00:52:11"
$wsZhCn.Range("K4").Value = "2016-08-21 00:52:29"

# de-de sheet: Correspond Handback DateTime for fe721b5d-0fa7-434f-8789-e65e1e9fafb1 row
$wsDeDe.Range("K4").Value = "2016-08-21 00:52:35"
